# "Grupo de Mercancias listo 100%"
# The merchandise-group sheet is finished: tidy up its name (drop the accented,
# spaced name in favour of an underscored one), make it the active tab, and
# leave the selection where the author left off (cell B19).

$wb = $excel.ActiveWorkbook

# Sheet 3 is "Grupo de mercancía" (Rastros(TIF), Rastros(Rastro), Grupo de
# mercancía, Sitio de inspección).
$ws = $wb.Worksheets.Item(3)

# Rename it - this also repoints the sheet-scoped _xlnm._FilterDatabase
# defined name automatically.
$ws.Name = "Grupo_mercancia"

# Make it the active/selected sheet (was "Rastros(Rastro)").
$ws.Activate()

# Leave the cursor on B19, matching where work left off.
$ws.Range("B19").Select()
